$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Add($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$ws.Name = "13"

$ws.Cells.Item(1, 1).Value = "add(add(conditional(conditional(conditional(x, conditional(conditional(conditional(add(conditional(x, vel), conditional(vel, y)), x), conditional(add(add(x, conditional(vel, y)), conditional(x, y)), conditional(y, y))), vel)), add(add(vel, conditional(add(y, conditional(y, vel)), add(add(conditional(x, vel), add(x, y)), conditional(add(add(add(y, vel), y), conditional(y, y)), add(y, add(add(conditional(conditional(y, y), conditional(vel, vel)), add(add(x, y), add(vel, conditional(conditional(y, y), add(y, y))))), x)))))), conditional(conditional(conditional(conditional(y, y), x), conditional(add(y, conditional(add(x, x), y)), conditional(add(vel, x), conditional(add(x, conditional(x, x)), x)))), add(conditional(add(x, y), x), add(add(x, y), x))))), add(vel, vel)), conditional(x, add(x, conditional(x, vel)))), conditional(add(add(y, y), conditional(y, y)), add(add(x, add(x, y)), add(add(x, add(add(y, vel), add(x, y))), x))))"
$ws.Cells.Item(1, 2).Value = -1944.09
$ws.Cells.Item(1, 3).Value = -1983.23
$ws.Cells.Item(1, 4).Value = -1833.96
$ws.Cells.Item(1, 5).Value = -1818.42
$ws.Cells.Item(1, 6).Value = -1901.66
$ws.Cells.Item(1, 7).Value = -2186.24
$ws.Cells.Item(1, 8).Value = -2240.38
$ws.Cells.Item(1, 9).Value = -2185.09
$ws.Cells.Item(1, 10).Value = -1957.49
$ws.Cells.Item(1, 11).Value = -2034.13
$ws.Cells.Item(1, 12).Value = -1954.08
$ws.Cells.Item(1, 13).Value = -1882.37
$ws.Cells.Item(1, 14).Value = -1695.86
$ws.Cells.Item(1, 15).Value = -1783.76
$ws.Cells.Item(1, 16).Value = -1827.98
$ws.Cells.Item(1, 17).Value = -2003.01
$ws.Cells.Item(2, 1).Value = "conditional(add(add(y, add(y, y)), conditional(x, conditional(add(add(vel, conditional(y, conditional(x, conditional(add(x, vel), conditional(x, x))))), add(add(add(vel, x), add(y, add(y, y))), vel)), add(add(vel, vel), conditional(x, y))))), add(add(add(add(x, add(add(x, x), y)), add(y, x)), y), add(vel, x)))"
$ws.Cells.Item(2, 2).Value = -1532.59
$ws.Cells.Item(2, 3).Value = -1348.34
$ws.Cells.Item(2, 4).Value = -1454.23
$ws.Cells.Item(2, 5).Value = -1411.54
$ws.Cells.Item(2, 6).Value = -1712.81
$ws.Cells.Item(2, 7).Value = -2032.2
$ws.Cells.Item(2, 8).Value = -1951.16
$ws.Cells.Item(2, 9).Value = -1903.98
$ws.Cells.Item(2, 10).Value = -1652.67
$ws.Cells.Item(2, 11).Value = -1812.1
$ws.Cells.Item(2, 12).Value = -1937.5
$ws.Cells.Item(2, 13).Value = -2096.37
$ws.Cells.Item(2, 14).Value = -2144.27
$ws.Cells.Item(2, 15).Value = -2203.32
$ws.Cells.Item(2, 16).Value = -2265.96
$ws.Cells.Item(2, 17).Value = -2312.92
$ws.Cells.Item(3, 1).Value = "conditional(add(y, add(y, y)), add(add(add(add(vel, add(x, y)), conditional(x, conditional(add(x, vel), conditional(x, vel)))), add(add(x, add(x, x)), add(add(y, x), y))), conditional(x, add(conditional(y, add(x, add(add(add(conditional(vel, x), add(add(x, x), add(conditional(conditional(add(x, y), add(add(conditional(x, y), conditional(y, vel)), conditional(x, vel))), y), conditional(vel, y)))), add(x, add(add(add(x, x), add(add(x, y), x)), add(y, conditional(x, conditional(conditional(vel, add(conditional(vel, x), conditional(x, y))), conditional(vel, vel))))))), conditional(conditional(vel, y), x)))), x))))"
$ws.Cells.Item(3, 2).Value = -1367.22
$ws.Cells.Item(3, 3).Value = -1389.54
$ws.Cells.Item(3, 4).Value = -1204.39
$ws.Cells.Item(3, 5).Value = -1231.77
$ws.Cells.Item(3, 6).Value = -1499.2
$ws.Cells.Item(3, 7).Value = -1654.53
$ws.Cells.Item(3, 8).Value = -1507.57
$ws.Cells.Item(3, 9).Value = -1476.15
$ws.Cells.Item(3, 10).Value = -1589.78
$ws.Cells.Item(3, 11).Value = -1747.01
$ws.Cells.Item(3, 12).Value = -1898.52
$ws.Cells.Item(3, 13).Value = -1984.07
$ws.Cells.Item(3, 14).Value = -2039.43
$ws.Cells.Item(3, 15).Value = -2074.51
$ws.Cells.Item(3, 16).Value = -2145.55
$ws.Cells.Item(3, 17).Value = -2185.22
$ws.Cells.Item(4, 1).Value = "conditional(add(conditional(y, y), add(add(y, y), x)), add(conditional(vel, add(conditional(add(conditional(add(vel, vel), add(add(x, add(conditional(vel, add(x, x)), add(add(conditional(y, add(x, add(add(vel, x), conditional(x, x)))), conditional(vel, x)), add(y, x)))), x)), x), conditional(add(conditional(x, vel), add(x, y)), y)), vel)), add(add(vel, x), add(add(x, vel), add(add(add(conditional(x, conditional(x, x)), add(add(y, x), add(x, x))), add(add(conditional(conditional(add(x, conditional(y, vel)), y), x), x), add(add(add(add(conditional(x, x), conditional(add(x, y), conditional(y, vel))), add(x, x)), add(x, vel)), vel))), add(conditional(x, vel), add(x, y)))))))"
$ws.Cells.Item(4, 2).Value = -1405.15
$ws.Cells.Item(4, 3).Value = -1481.18
$ws.Cells.Item(4, 4).Value = -1363.35
$ws.Cells.Item(4, 5).Value = -1609.2
$ws.Cells.Item(4, 6).Value = -1718.3
$ws.Cells.Item(4, 7).Value = -1584.89
$ws.Cells.Item(4, 8).Value = -1701.55
$ws.Cells.Item(4, 9).Value = -1608.97
$ws.Cells.Item(4, 10).Value = -1320.43
$ws.Cells.Item(4, 11).Value = -1549.2
$ws.Cells.Item(4, 12).Value = -1644.83
$ws.Cells.Item(4, 13).Value = -1775.87
$ws.Cells.Item(4, 14).Value = -1870.54
$ws.Cells.Item(4, 15).Value = -1975.6
$ws.Cells.Item(4, 16).Value = -2054.64
$ws.Cells.Item(4, 17).Value = -2111.4
$ws.Cells.Item(5, 1).Value = "conditional(add(add(conditional(y, y), x), add(y, y)), add(add(add(add(add(add(x, x), add(add(add(vel, conditional(x, x)), add(x, add(x, x))), vel)), conditional(y, vel)), add(x, y)), add(vel, x)), add(x, add(conditional(add(add(add(y, x), add(x, x)), conditional(conditional(vel, x), add(x, add(x, x)))), y), add(conditional(y, vel), x)))))"
$ws.Cells.Item(5, 2).Value = -1640.62
$ws.Cells.Item(5, 3).Value = -1632.89
$ws.Cells.Item(5, 4).Value = -1551.49
$ws.Cells.Item(5, 5).Value = -1626.25
$ws.Cells.Item(5, 6).Value = -1739.95
$ws.Cells.Item(5, 7).Value = -1574.58
$ws.Cells.Item(5, 8).Value = -1713.68
$ws.Cells.Item(5, 9).Value = -1658.51
$ws.Cells.Item(5, 10).Value = -1275.59
$ws.Cells.Item(5, 11).Value = -1500.89
$ws.Cells.Item(5, 12).Value = -1686.05
$ws.Cells.Item(5, 13).Value = -1774.19
$ws.Cells.Item(5, 14).Value = -1857.16
$ws.Cells.Item(5, 15).Value = -1974.98
$ws.Cells.Item(5, 16).Value = -2039.42
$ws.Cells.Item(5, 17).Value = -2123.3
$ws.Cells.Item(6, 1).Value = "add(add(conditional(x, add(conditional(x, x), add(conditional(conditional(x, x), add(conditional(conditional(add(vel, y), conditional(x, x)), add(vel, conditional(y, add(add(y, vel), add(x, conditional(add(y, x), conditional(x, vel))))))), add(add(vel, vel), vel))), conditional(y, conditional(conditional(vel, y), add(y, y)))))), conditional(y, add(conditional(x, add(vel, x)), add(add(add(conditional(x, x), y), x), add(conditional(vel, vel), add(conditional(add(vel, y), add(vel, x)), x)))))), conditional(add(y, y), add(add(add(x, vel), x), add(add(y, add(y, add(x, x))), conditional(x, x)))))"
$ws.Cells.Item(6, 2).Value = -2265.91
$ws.Cells.Item(6, 3).Value = -2277.55
$ws.Cells.Item(6, 4).Value = -2235.43
$ws.Cells.Item(6, 5).Value = -2255.4
$ws.Cells.Item(6, 6).Value = -2272.83
$ws.Cells.Item(6, 7).Value = -2266.74
$ws.Cells.Item(6, 8).Value = -2221.49
$ws.Cells.Item(6, 9).Value = -2134.59
$ws.Cells.Item(6, 10).Value = -2062.02
$ws.Cells.Item(6, 11).Value = -1664.44
$ws.Cells.Item(6, 12).Value = -1653.9
$ws.Cells.Item(6, 13).Value = -1738.8
$ws.Cells.Item(6, 14).Value = -1834.98
$ws.Cells.Item(6, 15).Value = -1858.45
$ws.Cells.Item(6, 16).Value = -1917.04
$ws.Cells.Item(6, 17).Value = -2042.67
$ws.Cells.Item(7, 1).Value = "add(conditional(y, add(add(add(x, y), conditional(y, add(conditional(add(vel, y), add(conditional(x, y), conditional(vel, x))), conditional(conditional(add(add(add(add(y, x), vel), conditional(vel, conditional(add(x, y), vel))), add(add(conditional(conditional(y, vel), y), x), y)), x), conditional(conditional(x, x), add(conditional(add(add(conditional(y, vel), conditional(vel, conditional(conditional(conditional(conditional(y, x), conditional(y, y)), y), conditional(x, x)))), x), conditional(vel, vel)), conditional(y, conditional(y, conditional(vel, conditional(y, add(vel, conditional(conditional(add(x, x), add(x, y)), add(y, vel))))))))))))), add(y, x))), conditional(add(conditional(conditional(x, y), conditional(y, x)), vel), add(add(y, vel), y)))"
$ws.Cells.Item(7, 2).Value = -2200.93
$ws.Cells.Item(7, 3).Value = -1790.35
$ws.Cells.Item(7, 4).Value = -1943.59
$ws.Cells.Item(7, 5).Value = -2177.36
$ws.Cells.Item(7, 6).Value = -2425.9
$ws.Cells.Item(7, 7).Value = -2432.09
$ws.Cells.Item(7, 8).Value = -2394.43
$ws.Cells.Item(7, 9).Value = -2389.31
$ws.Cells.Item(7, 10).Value = -2260.74
$ws.Cells.Item(7, 11).Value = -2160.22
$ws.Cells.Item(7, 12).Value = -1992.69
$ws.Cells.Item(7, 13).Value = -1733.56
$ws.Cells.Item(7, 14).Value = -1753.67
$ws.Cells.Item(7, 15).Value = -1842.77
$ws.Cells.Item(7, 16).Value = -1840.54
$ws.Cells.Item(7, 17).Value = -2014.42
$ws.Cells.Item(8, 1).Value = "conditional(add(add(add(y, y), add(y, y)), conditional(y, x)), add(add(x, x), add(vel, add(add(add(add(x, x), conditional(x, add(add(y, add(conditional(x, x), x)), add(add(vel, conditional(x, y)), conditional(conditional(x, y), conditional(add(conditional(conditional(vel, y), add(x, x)), add(y, y)), add(add(y, conditional(y, y)), add(vel, vel)))))))), conditional(y, conditional(conditional(y, y), conditional(x, y)))), add(conditional(conditional(x, vel), add(add(x, vel), add(vel, conditional(conditional(x, conditional(y, vel)), conditional(add(add(y, y), add(y, vel)), conditional(add(conditional(add(vel, vel), conditional(conditional(conditional(vel, y), add(x, vel)), y)), y), add(vel, vel))))))), add(x, add(x, add(add(y, vel), add(y, add(conditional(y, y), add(x, add(x, y))))))))))))"
$ws.Cells.Item(8, 2).Value = -1300.48
$ws.Cells.Item(8, 3).Value = -1245.68
$ws.Cells.Item(8, 4).Value = -1176.57
$ws.Cells.Item(8, 5).Value = -1359.83
$ws.Cells.Item(8, 6).Value = -1500.25
$ws.Cells.Item(8, 7).Value = -1508.47
$ws.Cells.Item(8, 8).Value = -1525.02
$ws.Cells.Item(8, 9).Value = -1161.59
$ws.Cells.Item(8, 10).Value = -1792.99
$ws.Cells.Item(8, 11).Value = -1912.36
$ws.Cells.Item(8, 12).Value = -1947.28
$ws.Cells.Item(8, 13).Value = -2030.01
$ws.Cells.Item(8, 14).Value = -2064.98
$ws.Cells.Item(8, 15).Value = -2109.59
$ws.Cells.Item(8, 16).Value = -2162.34
$ws.Cells.Item(8, 17).Value = -2227.2
$ws.Cells.Item(9, 1).Value = "conditional(add(add(conditional(x, vel), y), y), add(add(x, add(y, add(conditional(x, add(vel, y)), add(add(x, x), add(conditional(add(x, conditional(y, conditional(x, x))), add(conditional(conditional(conditional(y, y), add(add(x, conditional(y, conditional(y, x))), add(conditional(y, add(conditional(x, x), add(x, y))), add(x, vel)))), x), x)), x))))), add(x, vel)))"
$ws.Cells.Item(9, 2).Value = -1969.85
$ws.Cells.Item(9, 3).Value = -2006.03
$ws.Cells.Item(9, 4).Value = -1960.89
$ws.Cells.Item(9, 5).Value = -1920.7
$ws.Cells.Item(9, 6).Value = -1844.64
$ws.Cells.Item(9, 7).Value = -1759.9
$ws.Cells.Item(9, 8).Value = -1630.09
$ws.Cells.Item(9, 9).Value = -1543.47
$ws.Cells.Item(9, 10).Value = -1227.27
$ws.Cells.Item(9, 11).Value = -1575.49
$ws.Cells.Item(9, 12).Value = -1722.7
$ws.Cells.Item(9, 13).Value = -1869.59
$ws.Cells.Item(9, 14).Value = -1928.89
$ws.Cells.Item(9, 15).Value = -2041.11
$ws.Cells.Item(9, 16).Value = -2082.8
$ws.Cells.Item(9, 17).Value = -2137.56
$ws.Cells.Item(10, 1).Value = "add(conditional(add(y, y), add(vel, add(x, add(add(add(x, add(x, x)), add(x, y)), y)))), conditional(conditional(x, x), conditional(add(x, x), add(add(conditional(add(conditional(vel, vel), y), conditional(vel, y)), vel), add(add(x, x), add(add(conditional(y, conditional(add(add(add(add(y, y), conditional(y, y)), x), vel), conditional(y, conditional(x, vel)))), vel), y))))))"
$ws.Cells.Item(10, 2).Value = -2468.96
$ws.Cells.Item(10, 3).Value = -2424.31
$ws.Cells.Item(10, 4).Value = -2469.91
$ws.Cells.Item(10, 5).Value = -2483.38
$ws.Cells.Item(10, 6).Value = -2469.63
$ws.Cells.Item(10, 7).Value = -2417.33
$ws.Cells.Item(10, 8).Value = -2343.88
$ws.Cells.Item(10, 9).Value = -2251.9
$ws.Cells.Item(10, 10).Value = -2015.02
$ws.Cells.Item(10, 11).Value = -2124.34
$ws.Cells.Item(10, 12).Value = -2045.81
$ws.Cells.Item(10, 13).Value = -1618.21
$ws.Cells.Item(10, 14).Value = -1725.54
$ws.Cells.Item(10, 15).Value = -1777.62
$ws.Cells.Item(10, 16).Value = -1801.76
$ws.Cells.Item(10, 17).Value = -1989.1

$ws.Activate()